$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variable_Eff")

# Fix variable efficiency bug: the relationship class column header should
# reference "unit__from_node" rather than "unit__to_node"
$ws.Range("C1").Value = "unit__from_node"

# Match the resulting selection left by the edit
$ws.Range("C7").Select()
